# Update the "N" summary rows in the "Control at T1 v. T2" sheet.
# For each listed row:
#   - Columns I:N change from 2 to 1
#   - Columns O:T change from the inline string "0 (P=0.050)" to the numeric value 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 9, 13, 19, 25, 31, 37, 42, 48, 52, 58, 64, 70, 76)

foreach ($r in $rows) {
    $ws.Range("I$r`:N$r").Value = 1
    $ws.Range("O$r`:T$r").Value = 0
}
